$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 1.02
$ws.Cells.Item(2, 3).Value = 1.032418752245418
$ws.Cells.Item(2, 4).Value = 1.035910719710084
$ws.Cells.Item(2, 5).Value = 1.031613825857425
$ws.Cells.Item(2, 6).Value = 1.038572821579439
$ws.Cells.Item(2, 9).Value = 1.034637508855029
$ws.Cells.Item(2, 10).Value = 1.037549079072951
$ws.Cells.Item(2, 11).Value = 1.038705970644532
$ws.Cells.Item(2, 12).Value = 1.034421436142875
$ws.Cells.Item(2, 13).Value = 1.041360471105125
$ws.Cells.Item(2, 14).Value = 1.039022517773801
$ws.Cells.Item(3, 2).Value = 1.02
$ws.Cells.Item(3, 3).Value = 1.034119722557666
$ws.Cells.Item(3, 4).Value = 1.037237620941556
$ws.Cells.Item(3, 5).Value = 1.033105564823505
$ws.Cells.Item(3, 6).Value = 1.041053180131144
$ws.Cells.Item(3, 9).Value = 1.035177153994262
$ws.Cells.Item(3, 10).Value = 1.03888852100973
$ws.Cells.Item(3, 11).Value = 1.039840990845459
$ws.Cells.Item(3, 12).Value = 1.035719945905549
$ws.Cells.Item(3, 13).Value = 1.043646466135196
$ws.Cells.Item(3, 14).Value = 1.040363861871743
$ws.Cells.Item(4, 2).Value = 1.02
$ws.Cells.Item(4, 3).Value = 1.03521631883731
$ws.Cells.Item(4, 4).Value = 1.038092617490407
$ws.Cells.Item(4, 5).Value = 1.034066739868556
$ws.Cells.Item(4, 6).Value = 1.042649756546444
$ws.Cells.Item(4, 9).Value = 1.035522662894668
$ws.Cells.Item(4, 10).Value = 1.03975093339028
$ws.Cells.Item(4, 11).Value = 1.040571301666554
$ws.Cells.Item(4, 12).Value = 1.036555603434578
$ws.Cells.Item(4, 13).Value = 1.045117020195439
$ws.Cells.Item(4, 14).Value = 1.041227498976793
$ws.Cells.Item(5, 2).Value = 1.02
$ws.Cells.Item(5, 3).Value = 1.035676378928873
$ws.Cells.Item(5, 4).Value = 1.038451211435967
$ws.Cells.Item(5, 5).Value = 1.034469858303432
$ws.Cells.Item(5, 6).Value = 1.043319000715439
$ws.Cells.Item(5, 9).Value = 1.035667043084421
$ws.Cells.Item(5, 10).Value = 1.040112479119369
$ws.Cells.Item(5, 11).Value = 1.040877351065374
$ws.Cells.Item(5, 12).Value = 1.036905837063567
$ws.Cells.Item(5, 13).Value = 1.045733220291134
$ws.Cells.Item(5, 14).Value = 1.041589558142288
$ws.Cells.Item(6, 2).Value = 1.02
$ws.Cells.Item(6, 3).Value = 1.035753569954511
$ws.Cells.Item(6, 4).Value = 1.038511371712291
$ws.Cells.Item(6, 5).Value = 1.034537487888423
$ws.Cells.Item(6, 6).Value = 1.043431256459071
$ws.Cells.Item(6, 9).Value = 1.035691234283002
$ws.Cells.Item(6, 10).Value = 1.040173125331259
$ws.Cells.Item(6, 11).Value = 1.040928681431003
$ws.Cells.Item(6, 12).Value = 1.036964580151192
$ws.Cells.Item(6, 13).Value = 1.045836565847445
$ws.Cells.Item(6, 14).Value = 1.041650290478755
$ws.Cells.Item(7, 2).Value = 1.02
$ws.Cells.Item(7, 3).Value = 1.035222469885265
$ws.Cells.Item(7, 4).Value = 1.038097412347023
$ws.Cells.Item(7, 5).Value = 1.034072130103927
$ws.Cells.Item(7, 6).Value = 1.04265870664104
$ws.Cells.Item(7, 9).Value = 1.035524595523489
$ws.Cells.Item(7, 10).Value = 1.039755768338407
$ws.Cells.Item(7, 11).Value = 1.040575394919145
$ws.Cells.Item(7, 12).Value = 1.036560287482713
$ws.Cells.Item(7, 13).Value = 1.04512526176901
$ws.Cells.Item(7, 14).Value = 1.041232340791101
$ws.Cells.Item(8, 2).Value = 1.02
$ws.Cells.Item(8, 3).Value = 1.032994452504694
$ws.Cells.Item(8, 4).Value = 1.036359905874205
$ws.Cells.Item(8, 5).Value = 1.032118821169921
$ws.Cells.Item(8, 6).Value = 1.039412836116818
$ws.Cells.Item(8, 9).Value = 1.034820651183743
$ws.Cells.Item(8, 10).Value = 1.038002649016134
$ws.Cells.Item(8, 11).Value = 1.039090418133385
$ws.Cells.Item(8, 12).Value = 1.034861228584751
$ws.Cells.Item(8, 13).Value = 1.042134850831017
$ws.Cells.Item(8, 14).Value = 1.039476731838305
$ws.Cells.Item(9, 2).Value = 1.02
$ws.Cells.Item(9, 3).Value = 1.029036538841209
$ws.Cells.Item(9, 4).Value = 1.033269985265333
$ws.Cells.Item(9, 5).Value = 1.02864481807635
$ws.Cells.Item(9, 6).Value = 1.033626734008461
$ws.Cells.Item(9, 9).Value = 1.033551651534316
$ws.Cells.Item(9, 10).Value = 1.03487979814529
$ws.Cells.Item(9, 11).Value = 1.03644149449646
$ws.Cells.Item(9, 12).Value = 1.031831577576594
$ws.Cells.Item(9, 13).Value = 1.036797073170843
$ws.Cells.Item(9, 14).Value = 1.036349446161029
$ws.Cells.Item(10, 2).Value = 1.02
$ws.Cells.Item(10, 3).Value = 1.026375265900729
$ws.Cells.Item(10, 4).Value = 1.031190135704435
$ws.Cells.Item(10, 5).Value = 1.026306190443726
$ws.Cells.Item(10, 6).Value = 1.029721431878959
$ws.Cells.Item(10, 9).Value = 1.0326858923409
$ws.Cells.Item(10, 10).Value = 1.032774259754895
$ws.Cells.Item(10, 11).Value = 1.034653007841186
$ws.Cells.Item(10, 12).Value = 1.029786785191292
$ws.Cells.Item(10, 13).Value = 1.033189614051325
$ws.Cells.Item(10, 14).Value = 1.034240917664611
$ws.Cells.Item(11, 2).Value = 1.02
$ws.Cells.Item(11, 3).Value = 1.025217250615146
$ws.Cells.Item(11, 4).Value = 1.030284606649805
$ws.Cells.Item(11, 5).Value = 1.025287926098294
$ws.Cells.Item(11, 6).Value = 1.028018305107158
$ws.Cells.Item(11, 9).Value = 1.03230619634258
$ws.Cells.Item(11, 10).Value = 1.0318567007089
$ws.Cells.Item(11, 11).Value = 1.033873028246897
$ws.Cells.Item(11, 12).Value = 1.028895200722633
$ws.Cells.Item(11, 13).Value = 1.03161525255074
$ws.Cells.Item(11, 14).Value = 1.033322055579525
$ws.Cells.Item(12, 2).Value = 1.02
$ws.Cells.Item(12, 3).Value = 1.024786236586676
$ws.Cells.Item(12, 4).Value = 1.029947492240568
$ws.Cells.Item(12, 5).Value = 1.0249088311574
$ws.Cells.Item(12, 6).Value = 1.027383806073543
$ws.Cells.Item(12, 9).Value = 1.032164425691444
$ws.Cells.Item(12, 10).Value = 1.031514979862898
$ws.Cells.Item(12, 11).Value = 1.033582457421904
$ws.Cells.Item(12, 12).Value = 1.028563078589181
$ws.Cells.Item(12, 13).Value = 1.031028555910913
$ws.Cells.Item(12, 14).Value = 1.032979849450725
$ws.Cells.Item(13, 2).Value = 1.02
$ws.Cells.Item(13, 3).Value = 1.024878730574731
$ws.Cells.Item(13, 4).Value = 1.03001983914234
$ws.Cells.Item(13, 5).Value = 1.024990187854898
$ws.Cells.Item(13, 6).Value = 1.027519994460433
$ws.Cells.Item(13, 9).Value = 1.032194869400958
$ws.Cells.Item(13, 10).Value = 1.031588321132547
$ws.Cells.Item(13, 11).Value = 1.033644824672622
$ws.Cells.Item(13, 12).Value = 1.028634363137805
$ws.Cells.Item(13, 13).Value = 1.031154491672924
$ws.Cells.Item(13, 14).Value = 1.03305329487339
$ws.Cells.Item(14, 2).Value = 1.02
$ws.Cells.Item(14, 3).Value = 1.025181640864594
$ws.Cells.Item(14, 4).Value = 1.030256756290852
$ws.Cells.Item(14, 5).Value = 1.025256607788013
$ws.Cells.Item(14, 6).Value = 1.027965895980464
$ws.Cells.Item(14, 9).Value = 1.032294492587487
$ws.Cells.Item(14, 10).Value = 1.031828472405091
$ws.Cells.Item(14, 11).Value = 1.033849027047435
$ws.Cells.Item(14, 12).Value = 1.028867766859129
$ws.Cells.Item(14, 13).Value = 1.031566795309094
$ws.Cells.Item(14, 14).Value = 1.033293787188286
$ws.Cells.Item(15, 2).Value = 1.02
$ws.Cells.Item(15, 3).Value = 1.025368157134906
$ws.Cells.Item(15, 4).Value = 1.030402627424146
$ws.Cells.Item(15, 5).Value = 1.025420642472264
$ws.Cells.Item(15, 6).Value = 1.028240379335608
$ws.Cells.Item(15, 9).Value = 1.032355776039049
$ws.Cells.Item(15, 10).Value = 1.031976317854586
$ws.Cells.Item(15, 11).Value = 1.033974729489558
$ws.Cells.Item(15, 12).Value = 1.029011448373636
$ws.Cells.Item(15, 13).Value = 1.031820574464876
$ws.Cells.Item(15, 14).Value = 1.033441842595278
$ws.Cells.Item(16, 2).Value = 1.02
$ws.Cells.Item(16, 3).Value = 1.026451998112162
$ws.Cells.Item(16, 4).Value = 1.031250127036284
$ws.Cells.Item(16, 5).Value = 1.026373649001667
$ws.Cells.Item(16, 6).Value = 1.029834202382899
$ws.Cells.Item(16, 9).Value = 1.032710989098741
$ws.Cells.Item(16, 10).Value = 1.032835030227475
$ws.Cells.Item(16, 11).Value = 1.034704654043722
$ws.Cells.Item(16, 12).Value = 1.029845824888407
$ws.Cells.Item(16, 13).Value = 1.033293834870873
$ws.Cells.Item(16, 14).Value = 1.034301774438231
$ws.Cells.Item(17, 2).Value = 1.02
$ws.Cells.Item(17, 3).Value = 1.027130328824234
$ws.Cells.Item(17, 4).Value = 1.03178040559114
$ws.Cells.Item(17, 5).Value = 1.026969924333582
$ws.Cells.Item(17, 6).Value = 1.030830680337202
$ws.Cells.Item(17, 9).Value = 1.032932507624013
$ws.Cells.Item(17, 10).Value = 1.033372098300688
$ws.Cells.Item(17, 11).Value = 1.035161017534101
$ws.Cells.Item(17, 12).Value = 1.030367539679188
$ws.Cells.Item(17, 13).Value = 1.034214635376307
$ws.Cells.Item(17, 14).Value = 1.034839605209711
$ws.Cells.Item(18, 2).Value = 1.02
$ws.Cells.Item(18, 3).Value = 1.02752544257084
$ws.Cells.Item(18, 4).Value = 1.03208923267224
$ws.Cells.Item(18, 5).Value = 1.027317180557376
$ws.Cells.Item(18, 6).Value = 1.031410745433321
$ws.Cells.Item(18, 9).Value = 1.033061251659171
$ws.Cells.Item(18, 10).Value = 1.033684798242885
$ws.Cells.Item(18, 11).Value = 1.035426672035844
$ws.Cells.Item(18, 12).Value = 1.030671252582067
$ws.Cells.Item(18, 13).Value = 1.034750539284219
$ws.Cells.Item(18, 14).Value = 1.035152749221693
$ws.Cells.Item(19, 2).Value = 1.02
$ws.Cells.Item(19, 3).Value = 1.027660074420822
$ws.Cells.Item(19, 4).Value = 1.032194454635729
$ws.Cells.Item(19, 5).Value = 1.027435494820516
$ws.Cells.Item(19, 6).Value = 1.031608337162708
$ws.Cells.Item(19, 9).Value = 1.033105071680562
$ws.Cells.Item(19, 10).Value = 1.03379132598896
$ws.Cells.Item(19, 11).Value = 1.035517163199889
$ws.Cells.Item(19, 12).Value = 1.030774710598194
$ws.Cells.Item(19, 13).Value = 1.034933069718379
$ws.Cells.Item(19, 14).Value = 1.035259428249386
$ws.Cells.Item(20, 2).Value = 1.02
$ws.Cells.Item(20, 3).Value = 1.02705760691338
$ws.Cells.Item(20, 4).Value = 1.031723561028518
$ws.Cells.Item(20, 5).Value = 1.026906005755054
$ws.Cells.Item(20, 6).Value = 1.03072388849718
$ws.Cells.Item(20, 9).Value = 1.032908788852756
$ws.Cells.Item(20, 10).Value = 1.03331453433473
$ws.Cells.Item(20, 11).Value = 1.035112109471961
$ws.Cells.Item(20, 12).Value = 1.03031162625706
$ws.Cells.Item(20, 13).Value = 1.034115965033166
$ws.Cells.Item(20, 14).Value = 1.034781959496319
$ws.Cells.Item(21, 2).Value = 1.02
$ws.Cells.Item(21, 3).Value = 1.025092465701372
$ws.Cells.Item(21, 4).Value = 1.030187011240801
$ws.Cells.Item(21, 5).Value = 1.025178177843714
$ws.Cells.Item(21, 6).Value = 1.027834641561373
$ws.Cells.Item(21, 9).Value = 1.032265176411581
$ws.Cells.Item(21, 10).Value = 1.031757778808791
$ws.Cells.Item(21, 11).Value = 1.033788918189178
$ws.Cells.Item(21, 12).Value = 1.028799061612152
$ws.Cells.Item(21, 13).Value = 1.031445435338077
$ws.Cells.Item(21, 14).Value = 1.033222993198971
$ws.Cells.Item(22, 2).Value = 1.02
$ws.Cells.Item(22, 3).Value = 1.023851825231024
$ws.Cells.Item(22, 4).Value = 1.029216511492151
$ws.Cells.Item(22, 5).Value = 1.024086801297737
$ws.Cells.Item(22, 6).Value = 1.026007136823883
$ws.Cells.Item(22, 9).Value = 1.031856255649654
$ws.Cells.Item(22, 10).Value = 1.030773775817996
$ws.Cells.Item(22, 11).Value = 1.032952038953504
$ws.Cells.Item(22, 12).Value = 1.027842557365968
$ws.Cells.Item(22, 13).Value = 1.029755293865902
$ws.Cells.Item(22, 14).Value = 1.032237592811063
$ws.Cells.Item(23, 2).Value = 1.02
$ws.Cells.Item(23, 3).Value = 1.024510001771129
$ws.Cells.Item(23, 4).Value = 1.029731415987634
$ws.Cells.Item(23, 5).Value = 1.024665843971196
$ws.Cells.Item(23, 6).Value = 1.026976988271895
$ws.Cells.Item(23, 9).Value = 1.032073439442087
$ws.Cells.Item(23, 10).Value = 1.031295914975089
$ws.Cells.Item(23, 11).Value = 1.033396158263263
$ws.Cells.Item(23, 12).Value = 1.028350146006833
$ws.Cells.Item(23, 13).Value = 1.030652339691235
$ws.Cells.Item(23, 14).Value = 1.032760473465649
$ws.Cells.Item(24, 2).Value = 1.02
$ws.Cells.Item(24, 3).Value = 1.027090468498194
$ws.Cells.Item(24, 4).Value = 1.031749248108091
$ws.Cells.Item(24, 5).Value = 1.026934889480687
$ws.Cells.Item(24, 6).Value = 1.030772146725429
$ws.Cells.Item(24, 9).Value = 1.032919507777806
$ws.Cells.Item(24, 10).Value = 1.033340546751513
$ws.Cells.Item(24, 11).Value = 1.035134210571414
$ws.Cells.Item(24, 12).Value = 1.030336892962164
$ws.Cells.Item(24, 13).Value = 1.034160553567961
$ws.Cells.Item(24, 14).Value = 1.034808008853718
$ws.Cells.Item(25, 2).Value = 1.02
$ws.Cells.Item(25, 3).Value = 1.030063656907836
$ws.Cells.Item(25, 4).Value = 1.034072242343767
$ws.Cells.Item(25, 5).Value = 1.02954683778813
$ws.Cells.Item(25, 6).Value = 1.035130783362377
$ws.Cells.Item(25, 9).Value = 1.033883159609586
$ws.Cells.Item(25, 10).Value = 1.035691220050267
$ws.Cells.Item(25, 11).Value = 1.037130211396782
$ws.Cells.Item(25, 12).Value = 1.032619149536559
$ws.Cells.Item(25, 13).Value = 1.038185417289949
$ws.Cells.Item(25, 14).Value = 1.037162020378182
